$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# Weekly titration reading for the blue tanks (new row 42)
$ws.Range("A42").Value = 20210615
$ws.Range("B42").Value = 2228.7220000000002
$ws.Range("C42").Value = 2224.4699999999998
$ws.Range("D42").Formula = "=100*(B42-C42)/C42"
$ws.Range("E42").Value = 180
$ws.Range("F42").Value = "CRM OPENED 20210526"

# Move the active selection down to the next empty row, as Excel does
# after typing a new row of data
$ws.Range("A43").Select()
